$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# Selection change
$ws.Range("H4:H21").Select()
$ws.Application.ActiveCell = $ws.Range("H4")

# H2 changes condition
$ws.Range("H2").Formula = "=IF(F2=patience-1,TRUE)"

# Rows 3 to 21: update C, F, H, I formulas
for ($r = 3; $r -le 21; $r++) {
    $prev = $r - 1
    $ws.Range("C$r").Formula = "=IF(H$r=TRUE,B$r,IF(D$r>epsilon,B$r,C$prev))"
    $ws.Range("F$r").Formula = "=IF(E$r=TRUE,0,MOD(F$prev+1,5))"
    $ws.Range("H$r").Formula = "=IF(F$prev=patience-1,TRUE)"
    $ws.Range("I$r").Formula = "=IF(C$r<>C$prev,A$r,I$prev)"
}
